$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (GitHub Actions data refresh): new Price/Volume(1h)
# figures for every coin row, plus two pairs of rows whose rank flipped
# (EnergySwap/Decentraland and Quant/NEARProtocol swapped places).
#
# The Price column (D) stores numbers that look textual, e.g. thousand-dot
# separated "27.446.91" or plain decimals like "311.83". Excel's Value
# setter auto-coerces a plain-decimal string into a real number (losing the
# original text representation / introducing float rounding), so force
# Text format on each Price cell before writing its new value.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.446.91'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.857.53'
$ws.Range('E3').Value = '  +0.70%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.83'
$ws.Range('E5').Value = '  +0.85%  '
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4766'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3791'
$ws.Range('E8').Value = '  +3.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07304'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9299'
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.63'
$ws.Range('E11').Value = '  +4.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07781'
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.862.97'
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.446'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.548'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.06'
$ws.Range('E16').Value = '  +1.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.011'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008792'
$ws.Range('E18').Value = '  +1.50%  '
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.473.85'
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.59'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.090'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.68'
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.930'
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.80'
$ws.Range('E25').Value = '  +1.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.43'
$ws.Range('E26').Value = '  +1.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.002'
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.19'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.933'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08880'
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.330'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.202'
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7508'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.569'
$ws.Range('E34').Value = '  +1.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.699'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02036'
$ws.Range('E36').Value = '  +3.84%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.119'
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5543'
$ws.Range('E38').Value = '  +5.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05266'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.985'
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.009'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.545'
$ws.Range('E42').Value = '  +3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1514'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.65'
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4856'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.09'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.658'
$ws.Range('E48').Value = '  +3.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.25'
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06098'
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9115'
$ws.Range('E51').Value = '  +2.30%  '
